$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = "1C3CCCCB7FN1124UC"
$ws.Range("K6").Value = "Radio Off"
$ws.Range("N6").Value = "Night Mode"
$ws.Range("N7").Value = "Day Mode"
$ws.Range("O6").Value = "Phone Pickup Long Press"
$ws.Range("O7").Value = "VR Long Press"
